# Applies the "Add files via upload" laundry_data.xlsx edit:
#   - D47 was mis-typed as text "3484"; correct it to the genuine number 3484.
#   - Append six new attendance rows (48-53) for 董事長室 / 黃金昇 covering
#     2025-08-13 .. 2025-08-22, mirroring the existing row layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47: D47 was stored as text; fix it to a real number -------------
$ws.Range("D47").Value = 3484

# --- Rows 48-53: new rows, columns A..W -----------------------------------
# Column order: A date, B dept, C (blank), D empId, E (blank),
# F..T garment counts, U name, V cost center, W registrant.
# A leading "'" forces Excel to keep a value as literal text instead of
# auto-converting look-alike numbers/dates.
$newRows = @(
    @("'2025-08-13", "董事長室", "", 3484, "", 0, 0, 0, 0, 4, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "黃金昇", "董事長室", "samuel.huang"),
    @("'2025-08-14", "董事長室", "", 3484, "", 0, 1, 0, 0, 1, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, "黃金昇", "董事長室", "samuel.huang"),
    @("'2025-08-21", "董事長室", "", 3484, "", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 3, 2, 0, 0, 0, "黃金昇", "董事長室", "samuel.huang"),
    @("'2025-08-21", "董事長室", "", 3484, "", 0, 1, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "黃金昇", "董事長室", "samuel.huang"),
    @("'2025-08-21", "董事長室", "", 3484, "", 0, 0, 0, 0, 2, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, "黃金昇", "董事長室", "samuel.huang"),
    @("'2025-08-22", "董事長室", "", "'3484", "", 0, 0, 0, 0, 4, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, "黃金昇", "董事長室", "samuel.huang")
)

$startRow = 48
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $targetRow = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($targetRow, $j + 1).Value = $rowValues[$j]
    }
}
